$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-05 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-06 Saturday", 2) | Out-Null
$d.Content.Find.Execute("425×7=", $true, $false, $false, $false, $false, $true, 1, $false, "304×8=", 2) | Out-Null
$d.Content.Find.Execute("468×9=", $true, $false, $false, $false, $false, $true, 1, $false, "450×3=", 2) | Out-Null
$d.Content.Find.Execute("766×7=", $true, $false, $false, $false, $false, $true, 1, $false, "876×8=", 2) | Out-Null
$d.Content.Find.Execute("173×8=", $true, $false, $false, $false, $false, $true, 1, $false, "333×9=", 2) | Out-Null
$d.Content.Find.Execute("101×6=", $true, $false, $false, $false, $false, $true, 1, $false, "323×6=", 2) | Out-Null
$d.Content.Find.Execute("173×4=", $true, $false, $false, $false, $false, $true, 1, $false, "845×8=", 2) | Out-Null
$d.Content.Find.Execute("562×3=", $true, $false, $false, $false, $false, $true, 1, $false, "947×3=", 2) | Out-Null
$d.Content.Find.Execute("421×5=", $true, $false, $false, $false, $false, $true, 1, $false, "833×4=", 2) | Out-Null
$d.Content.Find.Execute("605×5=", $true, $false, $false, $false, $false, $true, 1, $false, "951×5=", 2) | Out-Null
$d.Content.Find.Execute("291×5=", $true, $false, $false, $false, $false, $true, 1, $false, "572×3=", 2) | Out-Null
$d.Content.Find.Execute("252×9=", $true, $false, $false, $false, $false, $true, 1, $false, "341×6=", 2) | Out-Null
$d.Content.Find.Execute("947×6=", $true, $false, $false, $false, $false, $true, 1, $false, "787×6=", 2) | Out-Null
$d.Content.Find.Execute("897×5=", $true, $false, $false, $false, $false, $true, 1, $false, "389×4=", 2) | Out-Null
$d.Content.Find.Execute("565×9=", $true, $false, $false, $false, $false, $true, 1, $false, "483×8=", 2) | Out-Null
$d.Content.Find.Execute("694×9=", $true, $false, $false, $false, $false, $true, 1, $false, "742×4=", 2) | Out-Null
$d.Content.Find.Execute("386×6=", $true, $false, $false, $false, $false, $true, 1, $false, "453×7=", 2) | Out-Null
$d.Content.Find.Execute("655×9=", $true, $false, $false, $false, $false, $true, 1, $false, "547×2=", 2) | Out-Null
$d.Content.Find.Execute("732×9=", $true, $false, $false, $false, $false, $true, 1, $false, "907×7=", 2) | Out-Null
$d.Content.Find.Execute("655×4=", $true, $false, $false, $false, $false, $true, 1, $false, "350×9=", 2) | Out-Null
$d.Content.Find.Execute("401×3=", $true, $false, $false, $false, $false, $true, 1, $false, "837×3=", 2) | Out-Null
$d.Content.Find.Execute("709×8=", $true, $false, $false, $false, $false, $true, 1, $false, "493×5=", 2) | Out-Null
$d.Content.Find.Execute("353×8=", $true, $false, $false, $false, $false, $true, 1, $false, "640×3=", 2) | Out-Null
$d.Content.Find.Execute("638×2=", $true, $false, $false, $false, $false, $true, 1, $false, "891×7=", 2) | Out-Null
$d.Content.Find.Execute("973×9=", $true, $false, $false, $false, $false, $true, 1, $false, "961×7=", 2) | Out-Null
$d.Content.Find.Execute("632×6=", $true, $false, $false, $false, $false, $true, 1, $false, "569×3=", 2) | Out-Null
